$wb = $excel.ActiveWorkbook

# --- Sheet "isa_template": update Version value from 1.0.2 to 1.0.3 ---
$meta = $wb.Worksheets.Item("isa_template")
$meta.Range("B4").Value = "1.0.3"

# --- Sheet "MS": update header names and data row values ---
$ms = $wb.Worksheets.Item("MS")

# Header row 1
$ms.Range("K1").Value = "Component [instrument]"
$ms.Range("T1").Value = "Output [Data]"

# Data row 2
$ms.Range("D2").Value = "https://bioregistry.io/NCIT:C17156"
$ms.Range("G2").Value = "https://www.ebi.ac.uk/ols4/ontologies/ms/classes/http%253A%252F%252Fpurl.obolibrary.org%252Fobo%252FMS_1000130"
$ms.Range("P2").Value = "https://www.ebi.ac.uk/ols4/ontologies/ms/classes/http%253A%252F%252Fpurl.obolibrary.org%252Fobo%252FMS_1000073"
